$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 57575
$ws.Range("B2").Value = "Alana Ribeiro"
$ws.Range("C2").Value = "Recursos Humanos"
$ws.Range("D2").Value = "Outros"
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 45102
$ws.Range("G2").Value = 9781.41

# Row 3
$ws.Range("A3").Value = 89935
$ws.Range("B3").Value = "Luiza Costa"
$ws.Range("C3").Value = "Atendimento ao Cliente"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 45078
$ws.Range("G3").Value = 4628.96

# Row 4
$ws.Range("A4").Value = 26666
$ws.Range("B4").Value = "Ana Beatriz Barbosa"
$ws.Range("C4").Value = "Engenharia"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 45093
$ws.Range("G4").Value = 9264.23

# Row 5
$ws.Range("A5").Value = 53734
$ws.Range("B5").Value = "Emanuel Gomes"
$ws.Range("C5").Value = "Vendas"
$ws.Range("D5").Value = "Outros"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45078
$ws.Range("G5").Value = 3448.28

# Row 6
$ws.Range("A6").Value = 35172
$ws.Range("B6").Value = "Ana Carolina Castro"
$ws.Range("C6").Value = "Jurídico"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45104
$ws.Range("G6").Value = 6277.26

# Row 7
$ws.Range("A7").Value = 99011
$ws.Range("B7").Value = "Sr. Thales Aragão"
$ws.Range("C7").Value = "Atendimento ao Cliente"
$ws.Range("D7").Value = "Doença"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45098
$ws.Range("G7").Value = 10660.38

# Row 8
$ws.Range("A8").Value = 32133
$ws.Range("B8").Value = "Marcos Vinicius Souza"
$ws.Range("C8").Value = "Financeiro"
$ws.Range("D8").Value = "Viagem de negócios"
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 45083
$ws.Range("G8").Value = 8648.360000000001

# Row 9
$ws.Range("A9").Value = 39758
$ws.Range("B9").Value = "André Almeida"
$ws.Range("C9").Value = "P&D"
$ws.Range("D9").Value = "Consulta médica"
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 45092
$ws.Range("G9").Value = 3235.94

# Row 10
$ws.Range("A10").Value = 73703
$ws.Range("B10").Value = "Gustavo Henrique Porto"
$ws.Range("C10").Value = "P&D"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 45089
$ws.Range("G10").Value = 11584.13

# Row 11
$ws.Range("A11").Value = 37616
$ws.Range("B11").Value = "Luna da Conceição"
$ws.Range("C11").Value = "Operações"
$ws.Range("D11").Value = "Doença"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 45102
$ws.Range("G11").Value = 9691.360000000001
